$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 8, pushing rows 8..29 down to 9..30
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the new data
$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(8, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(8, 4).Value = 45148
$ws.Cells.Item(8, 5).Value = 15
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100101
$ws.Cells.Item(8, 8).Value = "Berries"
$ws.Cells.Item(8, 9).Value = 100101007
$ws.Cells.Item(8, 10).Value = "Kiwi"
$ws.Cells.Item(8, 11).Value = "Hayward"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 300
$ws.Cells.Item(8, 14).Value = 22000
$ws.Cells.Item(8, 15).Value = 23000
$ws.Cells.Item(8, 16).Value = 22500
$ws.Cells.Item(8, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(8, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(8, 19).Value = 1250
$ws.Cells.Item(8, 20).Value = 18
